$d = $word.ActiveDocument
$nbsp = [char]160

# 1. Re-create the "smarthosting" bookmark so the runtime reassigns it id 0
#    (Word automatically assigns bookmark w:id values; deleting and re-adding
#    the only bookmark causes it to be renumbered starting from 0).
$bm = $d.Bookmarks.Item("smarthosting")
$bmRange = $bm.Range
$bm.Delete()
$d.Bookmarks.Add("smarthosting", $bmRange)

# 2. Remove the trailing " such as" (plus the following non-breaking space)
#    from the "Arguably the hardest part..." sentence, leaving just
#    "...obtained from exchanges".
$rng = $d.Content
$rng.Find.Execute(" such as$nbsp") | Out-Null
$rng.Delete()

# 3. Remove the CryptoBridge and HitBTC HYPERLINK fields entirely (field
#    codes, instruction text, and display text) using Field.Delete so the
#    whole field structure (fldChar begin/separate/end + instrText runs)
#    disappears rather than just the visible text.
$f2 = $d.Fields.Item(2)
$f2.Delete()
$f1 = $d.Fields.Item(1)
$f1.Delete()

# 4. Remove the now-orphaned ", " separator text (comma + non-breaking
#    space) that used to sit between the two hyperlink fields, without
#    touching the following "." run's distinct (bold) formatting.
$rng2 = $d.Content
$rng2.Find.Execute(",$nbsp") | Out-Null
$rng2.Delete()
